$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 73, shifting existing rows 73-118 down to 74-119.
$ws.Rows.Item(73).Insert()

# Fill in the new row 73 with the new record's data.
$ws.Cells.Item(73, 1).Value = 11
$ws.Cells.Item(73, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(73, 3).Value = "Bíobío"
$ws.Cells.Item(73, 4).Value = 44762
$ws.Cells.Item(73, 5).Value = 8
$ws.Cells.Item(73, 6).Value = 100112021
$ws.Cells.Item(73, 7).Value = "Ají"
$ws.Cells.Item(73, 8).Value = "Inferno"
$ws.Cells.Item(73, 9).Value = "Primera"
$ws.Cells.Item(73, 10).Value = 25
$ws.Cells.Item(73, 11).Value = 19000
$ws.Cells.Item(73, 12).Value = 21000
$ws.Cells.Item(73, 13).Value = 20200
$ws.Cells.Item(73, 14).Value = "`$/caja 15 kilos"
$ws.Cells.Item(73, 15).Value = "Provincia de Huasco"
$ws.Cells.Item(73, 16).Value = 1347
$ws.Cells.Item(73, 17).Value = 15
$ws.Cells.Item(73, 18).Value = "Hortaliza"
